# Update the limma output table (rows 2-11) with the recomputed
# gene statistics, re-sorted by ascending P.Value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("gene1",  -1.17938635672426,  14.22088115326463, -4.205265466015536, 0.000028465026707544, 0.00028465026707544, 2.189932207218739),
    @("gene5",  -0.6324493433217881, 16.3443560659832, -3.70888153777021,  0.0002198669051132755, 0.001099334525566378, -0.102360119616943),
    @("gene4",   0.3639954574941058, 18.3113081730598,  3.174407221009857, 0.001548245361016576, 0.005160817870055252, -2.29158336899609),
    @("gene2",  -0.5663683886699804, 16.02534944916673, -3.074341718212893, 0.002168084464122608, 0.005420211160306521, -2.164716618960036),
    @("gene3",   0.5584724117960018, 13.58519099544315,  1.911613288845917, 0.05621723165508651,  0.112434463310173,   -4.608396723735679),
    @("gene8",   0.2351166788199528, 16.97684981952542,  1.613410667212506, 0.1069772502976971,   0.1782954171628285,  -5.777154420298248),
    @("gene9",   0.1561028444463624, 16.8669522266167,   1.043096188344876, 0.2971609858261296,   0.3943164205625868,  -6.503116479686732),
    @("gene7",  -0.1658160195647369, 16.4440834168446,  -1.004359804931692, 0.3154531364500694,   0.3943164205625868,  -6.445890375171244),
    @("gene10", -0.235124282266133,  14.83914834046315, -0.8968133847476746, 0.370038751762738,   0.4111541686252645,  -6.097000591079622),
    @("gene6",  -0.02286980402421968, 13.61861869798013, -0.07805199423126882, 0.9378026266757549, 0.9378026266757549, -6.378719327334932)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $values = $data[$i]

    $ws.Cells.Item($row, 1).Value = $values[0]
    $ws.Cells.Item($row, 2).Value = $values[1]
    $ws.Cells.Item($row, 3).Value = $values[2]
    $ws.Cells.Item($row, 4).Value = $values[3]
    $ws.Cells.Item($row, 5).Value = $values[4]
    $ws.Cells.Item($row, 6).Value = $values[5]
    $ws.Cells.Item($row, 7).Value = $values[6]
}
